# Restructure ontology: remove mfd_hab1=Urban if mfd_areatype=Urban
#
# For every data row (2 through the last used row) where this applies:
#   - column F (habitat_typenumber): 1210 -> 2100
#   - column N (mfd_hab1): Urban -> Wastewater (takes over what was in O)
#   - column O (mfd_hab2): Wastewater -> Influent (takes over what was in P)
#   - column P (mfd_hab3): removed / cleared (its old value "Influent" moved into O)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# habitat_typenumber ("1210") is a numeric-looking code but must stay text
# (it was stored as inline string text before the edit too). Force the
# column to Text format first so re-entering the numeric string "2100"
# doesn't get silently promoted to a real Number by Excel's General format.
$ws.Range("F2:F" + $lastRow).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $fVal = $ws.Cells.Item($r, 6).Value2
    $nVal = $ws.Cells.Item($r, 14).Value2
    $oVal = $ws.Cells.Item($r, 15).Value2

    if ($fVal -eq "1210") {
        $ws.Cells.Item($r, 6).Value2 = "2100"
    }

    if ($nVal -eq "Urban") {
        $ws.Cells.Item($r, 14).Value2 = "Wastewater"
    }

    if ($oVal -eq "Wastewater") {
        $ws.Cells.Item($r, 15).Value2 = "Influent"
    }

    # mfd_hab3 (column P) no longer holds a value for this row - its old
    # content ("Influent") has moved up into column O above. Clear the cell
    # outright (rather than leaving an empty string) so it drops out of the
    # row entirely, same as the header (row 1, which keeps its P1 label)
    # being left untouched.
    $ws.Cells.Item($r, 16).ClearContents()
}
